# Add a new "Greece" sheet with test data, modeled on the existing
# "Croatia" sheet (same layout/styling), positioned right after it and
# made the active tab - mirroring how the source sheet was duplicated.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Duplicate the Croatia sheet, placing the copy immediately after it.
$croatia.Copy($null, $croatia)

$greece = $wb.Worksheets.Item($croatia.Index + 1)
$greece.Name = "Greece"

# Fill in Greece-specific test data.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3206"

# Croatia is no longer the tab in focus - clear its special selection
# back to a whole-sheet selection like the other non-active sheets.
$croatia.Activate()
$croatia.Cells.Select()

# Greece becomes the newly active/selected sheet with B4 highlighted.
$greece.Activate()
$greece.Range("B4").Select()
